# Fix issue with new user skills and skill comparison.
#
# 1. Insert a new leading column ("Unnamed: 0") holding a 0-based row index
#    for each user, shifting all existing columns (Name .. Added Skills)
#    one position to the right (A -> B, B -> C, ... T -> U).
# 2. Clear the stray "Added Skills" values (0) that were left behind for
#    the "AA" and "Test1" users.
# 3. Append a new user row ("A") with a full set of skill scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert new column A, shifting the rest of the sheet right ---
$ws.Range("A1").EntireColumn.Insert()

# Header for the new index column (match the header row's bold/bordered
# style used by the rest of row 1 - copy formatting from its neighbor)
$ws.Cells.Item(1,1).Value = "Unnamed: 0"
$ws.Cells.Item(1,2).Copy()
$ws.Cells.Item(1,1).PasteSpecial(-4122)

# 0-based row index for each existing data row (rows 2..6 after the shift)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(6,1).Value = 4

# --- 2. Clear the leftover "Added Skills" (0) values for AA and Test1 ---
# Column "Added Skills" is now column U (21) after the insert.
$ws.Cells.Item(4,21).ClearContents()
$ws.Cells.Item(5,21).ClearContents()

# --- 3. Add the new user "A" on row 7 ---
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "A"
$ws.Cells.Item(7,3).Value = 5
$ws.Cells.Item(7,4).Value = 2.5
$ws.Cells.Item(7,5).Value = 5
$ws.Cells.Item(7,6).Value = 2.5
$ws.Cells.Item(7,7).Value = 2.5
$ws.Cells.Item(7,8).Value = 2.5
$ws.Cells.Item(7,9).Value = 2.5
$ws.Cells.Item(7,10).Value = 2.5
$ws.Cells.Item(7,11).Value = 2.5
$ws.Cells.Item(7,12).Value = 5
$ws.Cells.Item(7,13).Value = 5
$ws.Cells.Item(7,14).Value = 2.5
$ws.Cells.Item(7,15).Value = 2.5
$ws.Cells.Item(7,16).Value = 2.5
$ws.Cells.Item(7,17).Value = 2.5
$ws.Cells.Item(7,18).Value = 2.5
$ws.Cells.Item(7,19).Value = 2.5
$ws.Cells.Item(7,20).Value = 2.5
